# Auto-generated Excel COM-interop edit script
# Updates Price (D) and Volume(1h) (E) columns for the cryptos list,
# and fixes two coin-name/link/value swaps (rows 40/41 and 50/51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.363.48'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '3.754.07'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''594.61'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '''169.21'
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").Value = '3.751.45'
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.525'
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").Value = '''6.50'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '''0.0000277'
$ws.Range("E13").Value = '  +7.12%  '
$ws.Range("D14").Value = '''36.57'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '4.383.32'
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("D16").Value = '3.757.33'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '''18.65'
$ws.Range("E17").Value = '  +1.91%  '
$ws.Range("D18").Value = '67.439.41'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").Value = '''7.21'
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '''10.50'
$ws.Range("E21").Value = '  -4.11%  '
$ws.Range("D22").Value = '''468.57'
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").Value = '''0.719'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("E24").Value = '  -7.21%  '
$ws.Range("D25").Value = '''83.88'
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("D26").Value = '''2.22'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''12.15'
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("D28").Value = '''10.40'
$ws.Range("E28").Value = '  +4.46%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '''2.90'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '3.904.28'
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").Value = '''7.65'
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").Value = '''30.55'
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").Value = '''2.23'
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").Value = '''9.12'
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("D36").Value = '3.717.12'
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").Value = '''3.80'
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '''5.85'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").Value = '''0.997'
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '''0.311'
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D45").Value = '''8.73'
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("D47").Value = '''45.85'
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").Value = '''398.92'
$ws.Range("E48").Value = '  -3.52%  '
$ws.Range("D49").Value = '''0.000270'
$ws.Range("E49").Value = '  -7.10%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''140.02'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = '''39.45'
$ws.Range("E51").Value = '  +4.35%  '
